# Delete the "proposal" worksheet and rename the remaining data sheet to "in".

$wb = $excel.ActiveWorkbook

# Turn off alerts so deleting the sheet doesn't prompt for confirmation.
$excel.DisplayAlerts = $false

$proposalSheet = $wb.Worksheets.Item("proposal")
$proposalSheet.Delete()

$dataSheet = $wb.Worksheets.Item("Sheet 1 - 20230615-a1r-nc-sessi")
$dataSheet.Name = "in"

$excel.DisplayAlerts = $true
